$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 9.014799999999994
$ws.Range("B6").Value = 6.725699999999998
$ws.Range("B7").Value = 5.519299999999999
$ws.Range("C7").Value = -13.93239999999999
$ws.Range("B8").Value = 7.635500000000002
$ws.Range("C11").Value = -11.78070000000001
$ws.Range("C12").Value = -11.2167
$ws.Range("D12").Value = -7.291699999999999
$ws.Range("D13").Value = -8.432600000000006
$ws.Range("D14").Value = -7.983799999999997
$ws.Range("C15").Value = -14.67429999999999
$ws.Range("B16").Value = 7.587799999999996
$ws.Range("D16").Value = -8.6235
$ws.Range("D19").Value = -7.889999999999999
$ws.Range("B20").Value = 9.164200000000005
$ws.Range("C20").Value = -12.4567
$ws.Range("D20").Value = -7.878899999999999
$ws.Range("B21").Value = 8.885000000000005
$ws.Range("C21").Value = -12.32800000000001
$ws.Range("C22").Value = -12.8274
$ws.Range("D22").Value = -8.1288
$ws.Range("C23").Value = -12.32630000000001
$ws.Range("B28").Value = 5.992800000000002
$ws.Range("B29").Value = 4.879700000000002
$ws.Range("C29").Value = -10.50050000000001
$ws.Range("B30").Value = 5.042999999999998
$ws.Range("B32").Value = 7.319799999999998
$ws.Range("C34").Value = -11.59790000000001
$ws.Range("D36").Value = -8.268299999999998
$ws.Range("B40").Value = 9.220999999999995
$ws.Range("C42").Value = -12.4473
$ws.Range("C43").Value = -13.60569999999999
$ws.Range("D43").Value = -8.369499999999999
$ws.Range("C44").Value = -13.79839999999999
$ws.Range("C45").Value = -13.62749999999999
$ws.Range("B46").Value = 6.143099999999996
$ws.Range("C46").Value = -13.3097
$ws.Range("D46").Value = -8.366399999999995
$ws.Range("C50").Value = -13.61849999999999
$ws.Range("D50").Value = -8.2372
$ws.Range("B51").Value = 5.982300000000002
$ws.Range("C51").Value = -12.4184
$ws.Range("B52").Value = 5.4605
$ws.Range("B57").Value = 6.143999999999998
$ws.Range("C57").Value = -13.67539999999999
$ws.Range("B59").Value = 5.784899999999999
$ws.Range("B62").Value = 6.783800000000003
$ws.Range("C65").Value = -12.9755
$ws.Range("B66").Value = 5.428500000000002
$ws.Range("C66").Value = -11.6637
$ws.Range("C67").Value = -11.5688
$ws.Range("B73").Value = 8.458799999999997
$ws.Range("B74").Value = 9.17199999999999
$ws.Range("D76").Value = -7.442299999999993
$ws.Range("B77").Value = 8.968600000000006
$ws.Range("C79").Value = -11.15540000000001
$ws.Range("C84").Value = -13.06
$ws.Range("C87").Value = -14.1714
$ws.Range("B92").Value = 5.727799999999994
$ws.Range("C92").Value = -11.3111
$ws.Range("D95").Value = -8.181399999999998
$ws.Range("C97").Value = -11.69540000000001
$ws.Range("D97").Value = -8.842899999999997
$ws.Range("D99").Value = -7.768799999999999
$ws.Range("B100").Value = 6.762299999999995
